$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("G16").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_14_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);`n"
$ws.Range("G17").Value = "wait(3);`nvalidate1;`nwait(2);`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_15_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G18").Value = "wait(3);`nvalidate1;`nwait(2);`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_16_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G19").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_17_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G20").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_18_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G21").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_19_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G22").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_29_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`nSwitchApp(NATIVE_APP);`nwait(2);`nCheckUITextContains(Google);`npress_Key(Back);"
$ws.Range("G23").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_31_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);`npress_Key(Back);"
$ws.Range("G24").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_32_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G25").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_35_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G26").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_37_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);`npress_Key(Back);"
$ws.Range("G27").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_45_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nUIAutomatorScreenshot(VT328-19);`nvalidate4;`npress_Key(Back);`nvalidate5;"
$ws.Range("G28").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_21_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"
$ws.Range("G29").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_23_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);`npress_Key(Back);`nClickUIButtonText(OK);"
$ws.Range("G31").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_43_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);`npress_Key(Back);"
$ws.Range("G32").Value = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT328_30_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`npress_Key(Back);"

$ws.Range("H27").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Intent JS Test`n};`nvalidate3`n{`nvalidate_OldText_Exists=VT328_45`n};`nvalidate4`n{`nvalidate_App_Launched_Device=com.android.gallery3d`n};`nvalidate5`n{`nvalidate_Screenshot=VT328-19`n};"

$ws.Range("A1").Select()